$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 11540703
$ws.Range("I40").Value = 1640.0769
$ws.Range("J40").Value = 23079766
$ws.Range("K40").Value = 1640.0769
$ws.Range("L40").Value = 23079766
$ws.Range("M40").Value = -1465.0769
$ws.Range("N40").Value = -23080116

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1999.8334
$ws.Range("I52").Value = 1999.8334
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 5999.5002
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -5839.5002
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3075148
$ws.Range("J112").Value = 4160203.8
$ws.Range("L112").Value = 12480611.4
$ws.Range("N112").Value = -12482827.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 38462460
$ws.Range("I118").Value = 38462460
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 115387380
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -115385723
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1351.9242
$ws.Range("I132").Value = 1251.4237
$ws.Range("K132").Value = 3754.2711
$ws.Range("M132").Value = -1224.2711

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3096.7778
$ws.Range("I135").Value = 3096.7778
$ws.Range("K135").Value = 27871.0002
$ws.Range("M135").Value = -25336.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6320.52
$ws.Range("I137").Value = 9898.538
$ws.Range("K137").Value = 29695.614
$ws.Range("M137").Value = -27145.614

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2636.2
$ws.Range("I138").Value = 1990.8276
$ws.Range("J138").Value = 2899.8027
$ws.Range("K138").Value = 5972.4828
$ws.Range("L138").Value = 8699.408100000001
$ws.Range("M138").Value = -832.4827999999998
$ws.Range("N138").Value = -18979.4081

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20929.514
$ws.Range("I32").Value = 19514.605
$ws.Range("K32").Value = 19514.605
$ws.Range("M32").Value = -19227.605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1614.25
$ws.Range("I45").Value = 1058.6364
$ws.Range("K45").Value = 1058.6364
$ws.Range("M45").Value = -681.6364000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4565.98
$ws.Range("I74").Value = 4655.3193
$ws.Range("K74").Value = 4655.3193
$ws.Range("M74").Value = -3781.3193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4565.98
$ws.Range("I77").Value = 4655.3193
$ws.Range("K77").Value = 23276.5965
$ws.Range("M77").Value = -18908.5965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6510.75
$ws.Range("I122").Value = 4519.8857
$ws.Range("K122").Value = 13559.6571
$ws.Range("M122").Value = -11109.6571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2904.246
$ws.Range("I132").Value = 2066.8125
$ws.Range("K132").Value = 6200.4375
$ws.Range("M132").Value = -3670.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4209.9546
$ws.Range("I20").Value = 3140.5
$ws.Range("K20").Value = 3140.5
$ws.Range("M20").Value = -2893.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1496.2258
$ws.Range("I86").Value = 1358.069
$ws.Range("K86").Value = 1358.069
$ws.Range("M86").Value = -235.069

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1496.2258
$ws.Range("I89").Value = 1358.069
$ws.Range("K89").Value = 6790.344999999999
$ws.Range("M89").Value = -1174.344999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4272.5
$ws.Range("I107").Value = 4297.5
$ws.Range("J107").Value = 4247.5
$ws.Range("K107").Value = 4297.5
$ws.Range("L107").Value = 4247.5
$ws.Range("M107").Value = -2377.5
$ws.Range("N107").Value = -8087.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3875.5
$ws.Range("I7").Value = 58.8
$ws.Range("J7").Value = 6601.7144
$ws.Range("K7").Value = 58.8
$ws.Range("L7").Value = 6601.7144
$ws.Range("M7").Value = 54.2
$ws.Range("N7").Value = -6827.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26339192
$ws.Range("I31").Value = 1895.2273
$ws.Range("K31").Value = 1895.2273
$ws.Range("M31").Value = -1600.2273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 26339192
$ws.Range("I34").Value = 1895.2273
$ws.Range("K34").Value = 1895.2273
$ws.Range("M34").Value = -1693.2273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 13999
$ws.Range("I69").Value = 13999
$ws.Range("K69").Value = 13999
$ws.Range("M69").Value = -13250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 13999
$ws.Range("I72").Value = 13999
$ws.Range("K72").Value = 41997
$ws.Range("M72").Value = -38253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 5550.433
$ws.Range("I107").Value = 647.1111
$ws.Range("K107").Value = 647.1111
$ws.Range("M107").Value = 1272.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 69326
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 69326
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 69326
$ws.Range("N133").Value = -74386
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 10750.889
$ws.Range("I56").Value = 10750.889
$ws.Range("K56").Value = 10750.889
$ws.Range("M56").Value = -10220.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1779.5555
$ws.Range("J113").Value = 1895.6086
$ws.Range("L113").Value = 5686.825800000001
$ws.Range("N113").Value = -10026.8258

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 6976.6665
$ws.Range("I120").Value = 4465
$ws.Range("J120").Value = 12000
$ws.Range("K120").Value = 13395
$ws.Range("L120").Value = 36000
$ws.Range("M120").Value = -8557
$ws.Range("N120").Value = -45676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 1802.9166
$ws.Range("I133").Value = 954.375
$ws.Range("J133").Value = 3500
$ws.Range("K133").Value = 2863.125
$ws.Range("L133").Value = 10500
$ws.Range("M133").Value = 2196.875
$ws.Range("N133").Value = -20620

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 16668381
$ws.Range("I137").Value = 1208.826
$ws.Range("J137").Value = 39218084
$ws.Range("K137").Value = 3626.478
$ws.Range("L137").Value = 117654252
$ws.Range("M137").Value = 1473.522
$ws.Range("N137").Value = -117664452

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2219.7273
$ws.Range("I139").Value = 1300.0769
$ws.Range("J139").Value = 3548.111
$ws.Range("K139").Value = 3900.2307
$ws.Range("L139").Value = 10644.333
$ws.Range("M139").Value = 1239.7693
$ws.Range("N139").Value = -20924.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84.5625
$ws.Range("I2").Value = 74.44444
$ws.Range("K2").Value = 74.44444
$ws.Range("M2").Value = 38.55556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 5630.1665
$ws.Range("I99").Value = 5630.1665
$ws.Range("K99").Value = 5630.1665
$ws.Range("M99").Value = -3384.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 28355.775
$ws.Range("I132").Value = 28215.29
$ws.Range("K132").Value = 84645.87
$ws.Range("M132").Value = -82115.87

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 33821.895
$ws.Range("J136").Value = 33821.895
$ws.Range("L136").Value = 101465.685
$ws.Range("N136").Value = -106565.685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2806.2144
$ws.Range("I46").Value = 1680.7273
$ws.Range("K46").Value = 1680.7273
$ws.Range("M46").Value = -1492.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1884.275
$ws.Range("I68").Value = 1780.8334
$ws.Range("J68").Value = 2815.25
$ws.Range("K68").Value = 1780.8334
$ws.Range("L68").Value = 2815.25
$ws.Range("M68").Value = -1031.8334
$ws.Range("N68").Value = -4313.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1884.275
$ws.Range("I71").Value = 1780.8334
$ws.Range("J71").Value = 2815.25
$ws.Range("K71").Value = 8904.166999999999
$ws.Range("L71").Value = 14076.25
$ws.Range("M71").Value = -5160.166999999999
$ws.Range("N71").Value = -21564.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 49998
$ws.Range("J96").Value = 49998
$ws.Range("L96").Value = 49998
$ws.Range("N96").Value = -55490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4128.875
$ws.Range("I100").Value = 4227.7144
$ws.Range("K100").Value = 4227.7144
$ws.Range("M100").Value = -3686.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 38750
$ws.Range("J125").Value = 38750
$ws.Range("L125").Value = 38750
$ws.Range("N125").Value = -48590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6790.961
$ws.Range("I132").Value = 2172.0454
$ws.Range("K132").Value = 6516.1362
$ws.Range("M132").Value = -3986.1362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2298.8918
$ws.Range("I136").Value = 1687.3103
$ws.Range("J136").Value = 4515.875
$ws.Range("K136").Value = 5061.9309
$ws.Range("L136").Value = 13547.625
$ws.Range("M136").Value = -2511.9309
$ws.Range("N136").Value = -18647.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3305.182
$ws.Range("I126").Value = 1978
$ws.Range("K126").Value = 5934
$ws.Range("M126").Value = -3464

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13335326
$ws.Range("I132").Value = 30304260
$ws.Range("J132").Value = 2592.4524
$ws.Range("K132").Value = 90912780
$ws.Range("L132").Value = 7777.3572
$ws.Range("M132").Value = -90910250
$ws.Range("N132").Value = -12837.3572
